$d = $word.ActiveDocument

# Locate the "Referência do depoimento:" run and collapse the range to its end
# (right before the hyperlink that follows it).
$range = $d.Content
$found = $range.Find.Execute("Referência do depoimento:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$range.Collapse(0)  # wdCollapseEnd

# Track the insertion so it is recorded as its own revision/run, then accept
# just that revision. This keeps the inserted space as a distinct <w:r>
# (with its own, inherited formatting) instead of being silently coalesced
# into the neighbouring run during save.
$wasTracking = $d.TrackRevisions
$d.TrackRevisions = $true
$range.InsertAfter(" ")
$d.TrackRevisions = $wasTracking

if ($d.Revisions.Count -gt 0) {
    $d.Revisions.Item(1).Accept()
}
